$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Round row 5 values (B5:AH5) to 2 decimal places (custom accuracy),
#    matching the exact literal results of the source edit.
$ws.Range("B5").Value  = 24.5
$ws.Range("C5").Value  = 18.06
$ws.Range("D5").Value  = 1.38
$ws.Range("E5").Value  = 53.22
$ws.Range("F5").Value  = 43.66
$ws.Range("G5").Value  = 19.28
$ws.Range("H5").Value  = 69.31
$ws.Range("I5").Value  = 29.67
$ws.Range("J5").Value  = 13.16
$ws.Range("K5").Value  = 19.54
$ws.Range("L5").Value  = 21.37
$ws.Range("M5").Value  = 22.5
$ws.Range("N5").Value  = 6.16
$ws.Range("O5").Value  = 19.17
$ws.Range("P5").Value  = 27.26
$ws.Range("Q5").Value  = 16.14
$ws.Range("R5").Value  = 0.88
$ws.Range("S5").Value  = 0.95
$ws.Range("T5").Value  = 284.52
$ws.Range("U5").Value  = 53.48
$ws.Range("V5").Value  = 17.7
$ws.Range("W5").Value  = 35.97
$ws.Range("X5").Value  = 18.9
$ws.Range("Y5").Value  = 2.81
$ws.Range("Z5").Value  = 34.45
$ws.Range("AA5").Value = 15.63
$ws.Range("AB5").Value = 13.86
$ws.Range("AC5").Value = 16.3
$ws.Range("AD5").Value = 22.37
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 62.63
$ws.Range("AG5").Value = 9.98
$ws.Range("AH5").Value = 22.13

# 2. Remove row 6 entirely (data trimmed from 1000 -> fewer rows),
#    this also updates the sheet dimension from A1:AH6 to A1:AH5.
$ws.Rows.Item(6).Delete()

# 3. Narrow columns J (10) and AB (28) from width 8 to width 7.
#    ColumnWidth character units map to the raw OOXML "width" with a
#    fixed +5/6 offset, so subtract it to land exactly on 7.
$narrowWidth = 7 - 0.8333333333333334
$ws.Columns.Item(10).ColumnWidth = $narrowWidth
$ws.Columns.Item(28).ColumnWidth = $narrowWidth
